# Add a "users" column to the "project hours" sheet, listing the users
# associated with each project row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New header for column E
$ws.Cells.Item(1, 5).Value = "users"
# Match the header style used by the other header cells (B1:D1)
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$users = @(
    "['Chenghao DUAN', 'Arun Lakshmanan']",
    "['Sierra Young', 'Karun Koppula']",
    "['Berk Cagilci', 'Olivas Hernandez, Daniel', 'Daniel Olivas Hernandez']",
    "['Jiyang Chen']",
    "['Alexander Hill']",
    "['Kyle Pieper']"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $users[$i]
}
